$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 11 values
$ws.Cells.Item(2, 4).Value = 44204
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 110
$ws.Cells.Item(2, 14).Value = 7000
$ws.Cells.Item(2, 15).Value = 7500
$ws.Cells.Item(2, 16).Value = 7318
$ws.Cells.Item(2, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(2, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(2, 19).Value = 1045
$ws.Cells.Item(2, 20).Value = 7

# Row 3 <- old row 2 values
$ws.Cells.Item(3, 4).Value = 44553
$ws.Cells.Item(3, 12).Value = "Especial"
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 22000
$ws.Cells.Item(3, 15).Value = 22000
$ws.Cells.Item(3, 16).Value = 22000
$ws.Cells.Item(3, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(3, 19).Value = 3667
$ws.Cells.Item(3, 20).Value = 6

# Row 4 <- old row 3 values
$ws.Cells.Item(4, 4).Value = 44553
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 150
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 18000
$ws.Cells.Item(4, 16).Value = 18000
$ws.Cells.Item(4, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(4, 19).Value = 3000
$ws.Cells.Item(4, 20).Value = 6

# Row 5 <- old row 12 values
$ws.Cells.Item(5, 4).Value = 44558
$ws.Cells.Item(5, 12).Value = "Especial"
$ws.Cells.Item(5, 13).Value = 20
$ws.Cells.Item(5, 14).Value = 22000
$ws.Cells.Item(5, 15).Value = 22000
$ws.Cells.Item(5, 16).Value = 22000
$ws.Cells.Item(5, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(5, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(5, 19).Value = 3667
$ws.Cells.Item(5, 20).Value = 6

# Row 6 <- old row 13 values
$ws.Cells.Item(6, 4).Value = 44558
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 25
$ws.Cells.Item(6, 14).Value = 18000
$ws.Cells.Item(6, 15).Value = 18000
$ws.Cells.Item(6, 16).Value = 18000
$ws.Cells.Item(6, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(6, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(6, 19).Value = 3000
$ws.Cells.Item(6, 20).Value = 6

# Row 7 <- old row 5 values
$ws.Cells.Item(7, 4).Value = 44550
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 60
$ws.Cells.Item(7, 14).Value = 24000
$ws.Cells.Item(7, 15).Value = 24000
$ws.Cells.Item(7, 16).Value = 24000
$ws.Cells.Item(7, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(7, 18).Value = "Región Metropolitana"
$ws.Cells.Item(7, 19).Value = 3429
$ws.Cells.Item(7, 20).Value = 7

# Row 8 <- old row 6 values
$ws.Cells.Item(8, 4).Value = 44572
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 65
$ws.Cells.Item(8, 14).Value = 20000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 20000
$ws.Cells.Item(8, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(8, 18).Value = "Región Metropolitana"
$ws.Cells.Item(8, 19).Value = 3333
$ws.Cells.Item(8, 20).Value = 6

# Row 11 <- old row 4 values
$ws.Cells.Item(11, 4).Value = 44561
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 200
$ws.Cells.Item(11, 14).Value = 18000
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 18000
$ws.Cells.Item(11, 17).Value = "`$/bandeja 6 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(11, 19).Value = 3000
$ws.Cells.Item(11, 20).Value = 6

# Row 12 <- old row 7 values
$ws.Cells.Item(12, 4).Value = 44189
$ws.Cells.Item(12, 12).Value = "Especial"
$ws.Cells.Item(12, 13).Value = 20
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 15000
$ws.Cells.Item(12, 16).Value = 15000
$ws.Cells.Item(12, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(12, 19).Value = 2143
$ws.Cells.Item(12, 20).Value = 7

# Row 13 <- old row 8 values
$ws.Cells.Item(13, 4).Value = 44189
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 30
$ws.Cells.Item(13, 14).Value = 13000
$ws.Cells.Item(13, 15).Value = 13000
$ws.Cells.Item(13, 16).Value = 13000
$ws.Cells.Item(13, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(13, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(13, 19).Value = 1857
$ws.Cells.Item(13, 20).Value = 7

